$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 195
$ws.Range("I31").Value = 195
$ws.Range("K31").Value = 585
$ws.Range("M31").Value = -355

$ws.Range("H43").Value = 2925.6667
$ws.Range("I43").Value = 5000
$ws.Range("J43").Value = 1888.5
$ws.Range("K43").Value = 5000
$ws.Range("L43").Value = 1888.5
$ws.Range("M43").Value = -4931
$ws.Range("N43").Value = -2026.5

$ws.Range("H64").Value = 4314.457
$ws.Range("I64").Value = 4084.158
$ws.Range("J64").Value = 4587.9375
$ws.Range("K64").Value = 4084.158
$ws.Range("L64").Value = 4587.9375
$ws.Range("M64").Value = -3836.158
$ws.Range("N64").Value = -5083.9375

$ws.Range("H67").Value = 4314.457
$ws.Range("I67").Value = 4084.158
$ws.Range("J67").Value = 4587.9375
$ws.Range("K67").Value = 4084.158
$ws.Range("L67").Value = 4587.9375
$ws.Range("M67").Value = -3226.158
$ws.Range("N67").Value = -6303.9375

$ws.Range("H74").Value = 5101
$ws.Range("I74").Value = 4535.5557
$ws.Range("J74").Value = 5666.4443
$ws.Range("K74").Value = 4535.5557
$ws.Range("L74").Value = 5666.4443
$ws.Range("M74").Value = -3599.5557
$ws.Range("N74").Value = -7538.4443

$ws.Range("H77").Value = 5101
$ws.Range("I77").Value = 4535.5557
$ws.Range("J77").Value = 5666.4443
$ws.Range("K77").Value = 22677.7785
$ws.Range("L77").Value = 28332.2215
$ws.Range("M77").Value = -17997.7785
$ws.Range("N77").Value = -37692.2215

$ws.Range("H116").Value = 3853.762
$ws.Range("I116").Value = 4158.5454
$ws.Range("J116").Value = 3518.5
$ws.Range("K116").Value = 4158.5454
$ws.Range("L116").Value = 3518.5
$ws.Range("M116").Value = -716.5454
$ws.Range("N116").Value = -10402.5

$ws.Range("H136").Value = 26989.75
$ws.Range("J136").Value = 26989.75
$ws.Range("L136").Value = 26989.75
$ws.Range("N136").Value = -37189.75

$ws.Range("H140").Value = 47963.332
$ws.Range("J140").Value = 47963.332
$ws.Range("L140").Value = 47963.332
$ws.Range("N140").Value = -58323.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 2648.125
$ws.Range("I35").Value = 1883.5714
$ws.Range("K35").Value = 1883.5714
$ws.Range("M35").Value = -1477.5714

$ws.Range("H63").Value = 1597
$ws.Range("I63").Value = 1597
$ws.Range("K63").Value = 1597
$ws.Range("M63").Value = -911

$ws.Range("H66").Value = 1597
$ws.Range("I66").Value = 1597
$ws.Range("K66").Value = 7985
$ws.Range("M66").Value = -4553

$ws.Range("H74").Value = 16943276
$ws.Range("I74").Value = 14306476
$ws.Range("K74").Value = 14306476
$ws.Range("M74").Value = -14305602

$ws.Range("H77").Value = 16943276
$ws.Range("I77").Value = 14306476
$ws.Range("K77").Value = 71532380
$ws.Range("M77").Value = -71528012

$ws.Range("H122").Value = 3045.0908
$ws.Range("I122").Value = 2571.9
$ws.Range("J122").Value = 7777
$ws.Range("K122").Value = 7715.700000000001
$ws.Range("L122").Value = 23331
$ws.Range("M122").Value = -5265.700000000001
$ws.Range("N122").Value = -28231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 3950
$ws.Range("I36").Value = 2900
$ws.Range("J36").Value = 5000
$ws.Range("K36").Value = 2900
$ws.Range("L36").Value = 5000
$ws.Range("M36").Value = -2366
$ws.Range("N36").Value = -6068

$ws.Range("H80").Value = 313.875
$ws.Range("I80").Value = 387.6
$ws.Range("J80").Value = 191
$ws.Range("K80").Value = 387.6
$ws.Range("L80").Value = 191
$ws.Range("M80").Value = 610.4
$ws.Range("N80").Value = -2187

$ws.Range("H83").Value = 313.875
$ws.Range("I83").Value = 387.6
$ws.Range("J83").Value = 191
$ws.Range("K83").Value = 1938
$ws.Range("L83").Value = 955
$ws.Range("M83").Value = 3054
$ws.Range("N83").Value = -10939

$ws.Range("H86").Value = 6412.4287
$ws.Range("I86").Value = 9955.714
$ws.Range("J86").Value = 2869.1428
$ws.Range("K86").Value = 9955.714
$ws.Range("L86").Value = 2869.1428
$ws.Range("M86").Value = -8832.714
$ws.Range("N86").Value = -5115.1428

$ws.Range("H89").Value = 6412.4287
$ws.Range("I89").Value = 9955.714
$ws.Range("J89").Value = 2869.1428
$ws.Range("K89").Value = 49778.57
$ws.Range("L89").Value = 14345.714
$ws.Range("M89").Value = -44162.57
$ws.Range("N89").Value = -25577.714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 47
$ws.Range("I7").Value = 64
$ws.Range("J7").Value = 32.833332
$ws.Range("K7").Value = 64
$ws.Range("L7").Value = 32.833332
$ws.Range("M7").Value = 49
$ws.Range("N7").Value = -258.833332

$ws.Range("H31").Value = 2540.9858
$ws.Range("I31").Value = 913.069
$ws.Range("J31").Value = 3665.024
$ws.Range("K31").Value = 913.069
$ws.Range("L31").Value = 3665.024
$ws.Range("M31").Value = -618.069
$ws.Range("N31").Value = -4255.023999999999

$ws.Range("H34").Value = 2540.9858
$ws.Range("I34").Value = 913.069
$ws.Range("J34").Value = 3665.024
$ws.Range("K34").Value = 913.069
$ws.Range("L34").Value = 3665.024
$ws.Range("M34").Value = -711.069
$ws.Range("N34").Value = -4069.024

$ws.Range("H99").Value = 41991.56
$ws.Range("I99").Value = 64369.938
$ws.Range("J99").Value = 2207.7778
$ws.Range("K99").Value = 64369.938
$ws.Range("L99").Value = 2207.7778
$ws.Range("M99").Value = -62871.938
$ws.Range("N99").Value = -5203.7778

$ws.Range("H126").Value = 41991.56
$ws.Range("I126").Value = 64369.938
$ws.Range("J126").Value = 2207.7778
$ws.Range("K126").Value = 193109.814
$ws.Range("L126").Value = 6623.3334
$ws.Range("M126").Value = -190639.814
$ws.Range("N126").Value = -11563.3334

$ws.Range("H134").Value = 16130259
$ws.Range("I134").Value = 19231550
$ws.Range("J134").Value = 3542.8
$ws.Range("K134").Value = 57694650
$ws.Range("L134").Value = 10628.4
$ws.Range("M134").Value = -57692115
$ws.Range("N134").Value = -15698.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 640.1818
$ws.Range("J36").Value = 1000
$ws.Range("L36").Value = 3000
$ws.Range("N36").Value = -3338

$ws.Range("H60").Value = 861
$ws.Range("I60").Value = 101.666664
$ws.Range("K60").Value = 304.999992
$ws.Range("M60").Value = -53.99999200000002

$ws.Range("H68").Value = 1297.7593
$ws.Range("I68").Value = 705.94116
$ws.Range("J68").Value = 1569.6757
$ws.Range("K68").Value = 2117.82348
$ws.Range("L68").Value = 4709.0271
$ws.Range("M68").Value = -1306.82348
$ws.Range("N68").Value = -6331.0271

$ws.Range("H71").Value = 1297.7593
$ws.Range("I71").Value = 705.94116
$ws.Range("J71").Value = 1569.6757
$ws.Range("K71").Value = 6353.47044
$ws.Range("L71").Value = 14127.0813
$ws.Range("M71").Value = -2297.47044
$ws.Range("N71").Value = -22239.0813

$ws.Range("H113").Value = 696.7059
$ws.Range("I113").Value = 676.4783
$ws.Range("J113").Value = 739
$ws.Range("K113").Value = 2029.4349
$ws.Range("L113").Value = 2217
$ws.Range("M113").Value = 140.5651
$ws.Range("N113").Value = -6557

$ws.Range("H122").Value = 1223.125
$ws.Range("I122").Value = 903.5
$ws.Range("J122").Value = 1329.6666
$ws.Range("K122").Value = 8131.5
$ws.Range("L122").Value = 11966.9994
$ws.Range("M122").Value = -5681.5
$ws.Range("N122").Value = -16866.9994

$ws.Range("H132").Value = 1705.5555
$ws.Range("I132").Value = 2487.2727
$ws.Range("J132").Value = 1168.125
$ws.Range("K132").Value = 22385.4543
$ws.Range("L132").Value = 10513.125
$ws.Range("M132").Value = -19855.4543
$ws.Range("N132").Value = -15573.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4081
$ws.Range("I80").Value = 4247.222
$ws.Range("J80").Value = 2959
$ws.Range("K80").Value = 4247.222
$ws.Range("L80").Value = 2959
$ws.Range("M80").Value = -3249.222
$ws.Range("N80").Value = -4955

$ws.Range("H83").Value = 4081
$ws.Range("I83").Value = 4247.222
$ws.Range("J83").Value = 2959
$ws.Range("K83").Value = 21236.11
$ws.Range("L83").Value = 14795
$ws.Range("M83").Value = -16244.11
$ws.Range("N83").Value = -24779

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5443.645
$ws.Range("I136").Value = 3233.389
$ws.Range("J136").Value = 8504
$ws.Range("K136").Value = 9700.167000000001
$ws.Range("L136").Value = 25512
$ws.Range("M136").Value = -7150.167000000001
$ws.Range("N136").Value = -30612

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1242.0646
$ws.Range("I126").Value = 1069.6428
$ws.Range("J126").Value = 2851.3333
$ws.Range("K126").Value = 3208.9284
$ws.Range("L126").Value = 8553.999899999999
$ws.Range("M126").Value = -738.9284000000002
$ws.Range("N126").Value = -13493.9999

$ws.Range("H136").Value = 19286984
$ws.Range("I136").Value = 27328702
$ws.Range("J136").Value = 690512.9
$ws.Range("K136").Value = 81986106
$ws.Range("L136").Value = 2071538.7
$ws.Range("M136").Value = -81983556
$ws.Range("N136").Value = -2076638.7
